$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet (Hoja1 -> meetings) ---
$ws.Name = "meetings"

# --- Header row (row 1): D1 "user" -> "type", add E1 "title" ---
$ws.Range("D1").Value = "type"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "title"

# --- Row 2: zoom meeting, juan ---
$ws.Range("D2").Value = "zoom"
$ws.Range("E2").Value = "juan"

# --- Row 3: google_meet link, pedro ---
$ws.Range("B3").ClearContents()
$ws.Range("A3").Value = "https://meet.google.com/url_to_your_metting"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://meet.google.com/url_to_your_metting")
$ws.Range("C3").Value = 0.50277777777777777
$ws.Range("D3").Value = "google_meet"
$ws.Range("E3").Value = "pedro"
$ws.Rows(3).RowHeight = 28.8

# --- Row 4: zoom meeting, sofia ---
$ws.Range("C4").Value = 0.54444444444444395
$ws.Range("D4").Value = "zoom"
$ws.Range("E4").Value = "sofia"

# --- Row 5 (new row): google_meet link, raul ---
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "https://meet.google.com/url_to_your_metting"
$ws.Hyperlinks.Add($ws.Range("A5"), "https://meet.google.com/url_to_your_metting")

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = 0.75277777777777777

$ws.Range("D5").Value = "google_meet"
$ws.Range("E5").Value = "raul"

# --- Column widths ---
$ws.Columns(1).ColumnWidth = 40.1
$ws.Columns(2).ColumnWidth = 23.498697916666668

# --- Selection ---
$ws.Range("C9").Select() | Out-Null
